# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") previously stored the pitcher's
# strikeout count ("Strike#") for the game. This regenerates that column
# from the newly-computed K values (s_vals) produced by the stats re-run,
# row by row, for every data row on Sheet1 (rows 2-88).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> newly computed K value.
$kValues = @{
    2 = 1
    3 = 2
    4 = 1
    5 = 2
    7 = 1
    8 = 0
    9 = 1
    11 = 1
    12 = 0
    13 = 2
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 2
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 3
    31 = 3
    32 = 2
    33 = 2
    34 = 2
    35 = 1
    36 = 1
    37 = 0
    38 = 0
    39 = 2
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 2
    45 = 0
    46 = 3
    47 = 2
    48 = 2
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 1
    55 = 0
    56 = 1
    57 = 0
    58 = 0
    59 = 1
    60 = 2
    61 = 0
    62 = 0
    63 = 2
    64 = 0
    65 = 3
    66 = 2
    67 = 3
    68 = 1
    69 = 2
    70 = 2
    71 = 2
    72 = 0
    73 = 0
    74 = 1
    75 = 2
    76 = 1
    77 = 1
    78 = 1
    79 = 0
    80 = 0
    82 = 1
    83 = 1
    84 = 1
    85 = 1
    86 = 1
    87 = 1
    88 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Count) K values in column G"
